# Duplicate the existing 20-row random-number block (A1:C20) into
# A21:C40, and move the active selection to E19 (was E13).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("A21").Value = 0.09558636215638416
$ws.Range("B21").Value = 0.43014567904895973
$ws.Range("C21").Value = 0.9347512369179227
$ws.Range("A22").Value = 0.584985451860591
$ws.Range("B22").Value = 0.8331573120413628
$ws.Range("C22").Value = 0.3092515357453779
$ws.Range("A23").Value = 0.2541150455955481
$ws.Range("B23").Value = 0.6916013844821515
$ws.Range("C23").Value = 0.7604975177012897
$ws.Range("A24").Value = 0.16805844861551922
$ws.Range("B24").Value = 0.22811306315321656
$ws.Range("C24").Value = 0.6916150262344297
$ws.Range("A25").Value = 0.6102277376579327
$ws.Range("B25").Value = 0.3936799543107754
$ws.Range("C25").Value = 0.43110966795358363
$ws.Range("A26").Value = 0.9496752109236136
$ws.Range("B26").Value = 0.9265061165754044
$ws.Range("C26").Value = 0.11818201202586232
$ws.Range("A27").Value = 0.9244597616721542
$ws.Range("B27").Value = 0.608885765529084
$ws.Range("C27").Value = 0.7025359667354973
$ws.Range("A28").Value = 0.424243874366491
$ws.Range("B28").Value = 0.0958339398547069
$ws.Range("C28").Value = 0.12564639228445462
$ws.Range("A29").Value = 0.6488870519114712
$ws.Range("B29").Value = 0.6774862704237027
$ws.Range("C29").Value = 0.21943385355532774
$ws.Range("A30").Value = 0.9598868507240403
$ws.Range("B30").Value = 0.7708443963815717
$ws.Range("C30").Value = 0.24509829978459619
$ws.Range("A31").Value = 0.2546894194395577
$ws.Range("B31").Value = 0.7212320073959982
$ws.Range("C31").Value = 0.540518110444642
$ws.Range("A32").Value = 0.23232590780878914
$ws.Range("B32").Value = 0.5697926788631499
$ws.Range("C32").Value = 0.9607478733232954
$ws.Range("A33").Value = 0.43591726258290364
$ws.Range("B33").Value = 0.4849624498452769
$ws.Range("C33").Value = 0.8174209894147595
$ws.Range("A34").Value = 0.9766020297508428
$ws.Range("B34").Value = 0.3795805588922665
$ws.Range("C34").Value = 0.7208807390008642
$ws.Range("A35").Value = 0.9906821714808408
$ws.Range("B35").Value = 0.8029578843331492
$ws.Range("C35").Value = 0.25989262000867286
$ws.Range("A36").Value = 0.25321974100165023
$ws.Range("B36").Value = 0.344814067347768
$ws.Range("C36").Value = 0.852821021292405
$ws.Range("A37").Value = 0.5598425257121624
$ws.Range("B37").Value = 0.906239040815767
$ws.Range("C37").Value = 0.9077485771761173
$ws.Range("A38").Value = 0.08412969353971411
$ws.Range("B38").Value = 0.7051007053474422
$ws.Range("C38").Value = 0.0031123204474163835
$ws.Range("A39").Value = 0.17190784462471487
$ws.Range("B39").Value = 0.218632976322497
$ws.Range("C39").Value = 0.04021753469479383
$ws.Range("A40").Value = 0.26502759355714556
$ws.Range("B40").Value = 0.6175304988987728
$ws.Range("C40").Value = 0.004897938238207611

$ws.Range("E19").Select()
